$wb = $excel.ActiveWorkbook

# Change the selection on Simulation1 from H14 to B1:E1
$ws1 = $wb.Worksheets.Item("Simulation1")
$ws1.Range("B1:E1").Select() | Out-Null

# Add the new worksheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "FuncSimulation"

# Header row
$newSheet.Range("A1").Value = "Model"
$newSheet.Range("B1").Value = "T1_tpr"
$newSheet.Range("C1").Value = "T1_fpr"
$newSheet.Range("D1").Value = "T2_tpr"
$newSheet.Range("E1").Value = "T2_fpr"

# Row labels: write row 3's label before row 2's so the shared-string
# table indices line up with the target workbook (23 = Polynomial Degree,
# 24 = Model Predictors)
$newSheet.Range("A3").Value = "Polynomial Degree"
$newSheet.Range("A2").Value = "Model Predictors"

# Row 2 numeric data
$newSheet.Range("B2").Value = 1
$newSheet.Range("C2").Value = 0.1417
$newSheet.Range("D2").Value = 1
$newSheet.Range("E2").Value = 0.0037

# Row 3 numeric data
$newSheet.Range("B3").Value = 0.7125
$newSheet.Range("C3").Value = 0.0251
$newSheet.Range("D3").Value = 0.9667
$newSheet.Range("E3").Value = 0.0007

$newSheet.Columns.Item(1).AutoFit()

$newSheet.Range("A3").Select() | Out-Null
